# "fixed GET STARTED link"
# The "GET STARTED" button on the first slide (row 2) pointed to a deep
# sub-page (https://www.pdx.edu/civil-environmental-engineering/student-orgs).
# Fix it so it links to the PSU homepage instead (https://www.pdx.edu/),
# applying Excel's built-in "Hyperlink" cell style to match the other
# hyperlinked cells in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$buttonLink = $ws.Range("E2")
$buttonLink.Value = "https://www.pdx.edu/"
$buttonLink.Style = "Hyperlink"

# Leave the selection where the author ended up after making the edit.
[void]$ws.Range("F8").Select()
